$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values (column D) are plain digit/decimal
# strings (e.g. "597.52"). If written as-is, Excel would auto-convert them
# from text to a numeric value, which would both change their stored type
# and silently round/alter values like "1.00" -> 1 or "0.0000138" -> 1.38E-5.
# To keep them as plain text (as they are in the source workbook), we
# temporarily mark these specific cells as Text-formatted before writing,
# then restore their normal (default) style afterwards so no stray
# number-format is left behind.
$textForcedCells = @(
    "D5",
    "D6",
    "D13",
    "D18",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D30",
    "D32",
    "D33",
    "D37",
    "D38",
    "D39",
    "D42",
    "D43",
    "D45",
    "D46",
    "D48",
    "D49"
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Refresh Price (D) / Volume(1h) (E) figures pulled on Fri May 24 2024.
$ws.Range("D2").Value = '67.530.35'
$ws.Range("E2").Value = '  -3.01%  '
$ws.Range("D3").Value = '3.713.39'
$ws.Range("E3").Value = '  -4.63%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '597.52'
$ws.Range("E5").Value = '  -1.70%  '
$ws.Range("D6").Value = '166.91'
$ws.Range("E6").Value = '  -3.86%  '
$ws.Range("D7").Value = '3.712.31'
$ws.Range("E7").Value = '  -3.47%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("E11").Value = '  -2.56%  '
$ws.Range("E12").Value = '  -3.54%  '
$ws.Range("D13").Value = '37.84'
$ws.Range("E13").Value = '  -4.60%  '
$ws.Range("E14").Value = '  -4.11%  '
$ws.Range("D15").Value = '4.331.02'
$ws.Range("E15").Value = '  -5.46%  '
$ws.Range("D16").Value = '3.710.39'
$ws.Range("E16").Value = '  -5.78%  '
$ws.Range("D17").Value = '67.518.86'
$ws.Range("E17").Value = '  -3.31%  '
$ws.Range("D18").Value = '7.27'
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("D19").Value = '17.64'
$ws.Range("E19").Value = '  +6.82%  '
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").Value = '487.15'
$ws.Range("E21").Value = '  -3.23%  '
$ws.Range("D22").Value = '9.28'
$ws.Range("E22").Value = '  -3.85%  '
$ws.Range("D23").Value = '0.728'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("D24").Value = '85.21'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  -5.45%  '
$ws.Range("D26").Value = '0.0000138'
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  -3.17%  '
$ws.Range("D28").Value = '10.10'
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").Value = '2.94'
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("E31").Value = '  -7.19%  '
$ws.Range("D32").Value = '7.72'
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("D33").Value = '31.37'
$ws.Range("E33").Value = '  -6.42%  '
$ws.Range("D34").Value = '3.849.43'
$ws.Range("E34").Value = '  -5.70%  '
$ws.Range("E35").Value = '  -4.03%  '
$ws.Range("D36").Value = '3.651.30'
$ws.Range("E36").Value = '  -5.12%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -3.00%  '
$ws.Range("D39").Value = '5.82'
$ws.Range("E39").Value = '  -4.28%  '
$ws.Range("E40").Value = '  -7.52%  '
$ws.Range("E41").Value = '  -3.10%  '
$ws.Range("D42").Value = '48.84'
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("D43").Value = '425.96'
$ws.Range("E43").Value = '  -10.27%  '
$ws.Range("E44").Value = '  -4.63%  '
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -3.83%  '
$ws.Range("D46").Value = '8.46'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = '40.37'
$ws.Range("E48").Value = '  -4.72%  '
$ws.Range("D49").Value = '140.28'
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("E50").Value = '  -3.63%  '
$ws.Range("D51").Value = '2.742.17'
$ws.Range("E51").Value = '  -6.39%  '

# Restore the default "Normal" cell style now that the text values are set,
# so the cells end up with no special number format applied.
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
